# Update NATMI ligand-receptor TPM-derived statistics (Adam9-Itgb1) with the
# new TPM values. Columns affected (G,H,I,J = ligand stats; M,N,O,P = receptor
# stats; Q,R,S,T = edge stats) for data rows 2-26. Columns A-F, K, L are left
# untouched because the underlying cluster/gene identifiers did not change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 9.462749333333333
$ws.Cells.Item(2, 8).Value = 28.388248
$ws.Cells.Item(2, 9).Value = 0.07254428564686972
$ws.Cells.Item(2, 10).Value = 0.07439525120506714
$ws.Cells.Item(2, 13).Value = 121.928739
$ws.Cells.Item(2, 14).Value = 365.786217
$ws.Cells.Item(2, 15).Value = 0.2282232151508951
$ws.Cells.Item(2, 16).Value = 0.2419720431319445
$ws.Cells.Item(2, 17).Value = 1153.781093686424
$ws.Cells.Item(2, 18).Value = 10384.02984317781
$ws.Cells.Item(2, 19).Value = 0.01655629011115354
$ws.Cells.Item(2, 20).Value = 0.01800157093340436

$ws.Cells.Item(3, 7).Value = 9.462749333333333
$ws.Cells.Item(3, 8).Value = 28.388248
$ws.Cells.Item(3, 9).Value = 0.07254428564686972
$ws.Cells.Item(3, 10).Value = 0.07439525120506714
$ws.Cells.Item(3, 13).Value = 147.91433
$ws.Cells.Item(3, 14).Value = 443.74299
$ws.Cells.Item(3, 15).Value = 0.2768624053389947
$ws.Cells.Item(3, 16).Value = 0.2935413991166814
$ws.Cells.Item(3, 17).Value = 1399.676227597946
$ws.Cells.Item(3, 18).Value = 12597.08604838152
$ws.Cells.Item(3, 19).Value = 0.02008478541779146
$ws.Cells.Item(3, 20).Value = 0.02183808612637239

$ws.Cells.Item(4, 7).Value = 9.462749333333333
$ws.Cells.Item(4, 8).Value = 28.388248
$ws.Cells.Item(4, 9).Value = 0.07254428564686972
$ws.Cells.Item(4, 10).Value = 0.07439525120506714
$ws.Cells.Item(4, 13).Value = 83.50496933333334
$ws.Cells.Item(4, 14).Value = 250.514908
$ws.Cells.Item(4, 15).Value = 0.1563025480180701
$ws.Cells.Item(4, 16).Value = 0.1657186665504434
$ws.Cells.Item(4, 17).Value = 790.1865928890204
$ws.Cells.Item(4, 18).Value = 7111.679336001183
$ws.Cells.Item(4, 19).Value = 0.01133885669075645
$ws.Cells.Item(4, 20).Value = 0.01232868182738899

$ws.Cells.Item(5, 7).Value = 9.462749333333333
$ws.Cells.Item(5, 8).Value = 28.388248
$ws.Cells.Item(5, 9).Value = 0.07254428564686972
$ws.Cells.Item(5, 10).Value = 0.07439525120506714
$ws.Cells.Item(5, 13).Value = 91.06846250000001
$ws.Cells.Item(5, 14).Value = 182.136925
$ws.Cells.Item(5, 15).Value = 0.1704597085236707
$ws.Cells.Item(5, 16).Value = 0.1204857969594293
$ws.Cells.Item(5, 17).Value = 861.7580328095668
$ws.Cells.Item(5, 18).Value = 5170.5481968574
$ws.Cells.Item(5, 19).Value = 0.01236587778642332
$ws.Cells.Item(5, 20).Value = 0.00896357113143946

$ws.Cells.Item(6, 7).Value = 9.462749333333333
$ws.Cells.Item(6, 8).Value = 28.388248
$ws.Cells.Item(6, 9).Value = 0.07254428564686972
$ws.Cells.Item(6, 10).Value = 0.07439525120506714
$ws.Cells.Item(6, 13).Value = 89.83562999999999
$ws.Cells.Item(6, 14).Value = 269.50689
$ws.Cells.Item(6, 15).Value = 0.1681521229683693
$ws.Cells.Item(6, 16).Value = 0.1782820942415013
$ws.Cells.Item(6, 17).Value = 850.0920478920799
$ws.Cells.Item(6, 18).Value = 7650.828431028719
$ws.Cells.Item(6, 19).Value = 0.01219847564074495
$ws.Cells.Item(6, 20).Value = 0.01326334118646194

$ws.Cells.Item(7, 7).Value = 42.14988333333334
$ws.Cells.Item(7, 8).Value = 126.44965
$ws.Cells.Item(7, 9).Value = 0.3231336970688258
$ws.Cells.Item(7, 10).Value = 0.3313784449305509
$ws.Cells.Item(7, 13).Value = 121.928739
$ws.Cells.Item(7, 14).Value = 365.786217
$ws.Cells.Item(7, 15).Value = 0.2282232151508951
$ws.Cells.Item(7, 16).Value = 0.2419720431319445
$ws.Cells.Item(7, 17).Value = 5139.28212383045
$ws.Cells.Item(7, 18).Value = 46253.53911447405
$ws.Cells.Item(7, 19).Value = 0.07374661126864279
$ws.Cells.Item(7, 20).Value = 0.08018431936973196

$ws.Cells.Item(8, 7).Value = 42.14988333333334
$ws.Cells.Item(8, 8).Value = 126.44965
$ws.Cells.Item(8, 9).Value = 0.3231336970688258
$ws.Cells.Item(8, 10).Value = 0.3313784449305509
$ws.Cells.Item(8, 13).Value = 147.91433
$ws.Cells.Item(8, 14).Value = 443.74299
$ws.Cells.Item(8, 15).Value = 0.2768624053389947
$ws.Cells.Item(8, 16).Value = 0.2935413991166814
$ws.Cells.Item(8, 17).Value = 6234.571752828166
$ws.Cells.Item(8, 18).Value = 56111.1457754535
$ws.Cells.Item(8, 19).Value = 0.08946357261655719
$ws.Cells.Item(8, 20).Value = 0.09727329236202406

$ws.Cells.Item(9, 7).Value = 42.14988333333334
$ws.Cells.Item(9, 8).Value = 126.44965
$ws.Cells.Item(9, 9).Value = 0.3231336970688258
$ws.Cells.Item(9, 10).Value = 0.3313784449305509
$ws.Cells.Item(9, 13).Value = 83.50496933333334
$ws.Cells.Item(9, 14).Value = 250.514908
$ws.Cells.Item(9, 15).Value = 0.1563025480180701
$ws.Cells.Item(9, 16).Value = 0.1657186665504434
$ws.Cells.Item(9, 17).Value = 3519.724715153578
$ws.Cells.Item(9, 18).Value = 31677.5224363822
$ws.Cells.Item(9, 19).Value = 0.05050662020235665
$ws.Cells.Item(9, 20).Value = 0.05491559401745042

$ws.Cells.Item(10, 7).Value = 42.14988333333334
$ws.Cells.Item(10, 8).Value = 126.44965
$ws.Cells.Item(10, 9).Value = 0.3231336970688258
$ws.Cells.Item(10, 10).Value = 0.3313784449305509
$ws.Cells.Item(10, 13).Value = 91.06846250000001
$ws.Cells.Item(10, 14).Value = 182.136925
$ws.Cells.Item(10, 15).Value = 0.1704597085236707
$ws.Cells.Item(10, 16).Value = 0.1204857969594293
$ws.Cells.Item(10, 17).Value = 3838.525069721042
$ws.Cells.Item(10, 18).Value = 23031.15041832626
$ws.Cells.Item(10, 19).Value = 0.05508127581652817
$ws.Cells.Item(10, 20).Value = 0.03992639603263379

$ws.Cells.Item(11, 7).Value = 42.14988333333334
$ws.Cells.Item(11, 8).Value = 126.44965
$ws.Cells.Item(11, 9).Value = 0.3231336970688258
$ws.Cells.Item(11, 10).Value = 0.3313784449305509
$ws.Cells.Item(11, 13).Value = 89.83562999999999
$ws.Cells.Item(11, 14).Value = 269.50689
$ws.Cells.Item(11, 15).Value = 0.1681521229683693
$ws.Cells.Item(11, 16).Value = 0.1782820942415013
$ws.Cells.Item(11, 17).Value = 3786.5613236765
$ws.Cells.Item(11, 18).Value = 34079.0519130885
$ws.Cells.Item(11, 19).Value = 0.05433561716474099
$ws.Cells.Item(11, 20).Value = 0.05907884314871061

$ws.Cells.Item(12, 7).Value = 32.300192
$ws.Cells.Item(12, 8).Value = 96.900576
$ws.Cells.Item(12, 9).Value = 0.2476229975407503
$ws.Cells.Item(12, 10).Value = 0.2539410918713864
$ws.Cells.Item(12, 13).Value = 121.928739
$ws.Cells.Item(12, 14).Value = 365.786217
$ws.Cells.Item(12, 15).Value = 0.2282232151508951
$ws.Cells.Item(12, 16).Value = 0.2419720431319445
$ws.Cells.Item(12, 17).Value = 3938.321680017888
$ws.Cells.Item(12, 18).Value = 35444.89512016099
$ws.Cells.Item(12, 19).Value = 0.05651331664405222
$ws.Cells.Item(12, 20).Value = 0.0614466448352762

$ws.Cells.Item(13, 7).Value = 32.300192
$ws.Cells.Item(13, 8).Value = 96.900576
$ws.Cells.Item(13, 9).Value = 0.2476229975407503
$ws.Cells.Item(13, 10).Value = 0.2539410918713864
$ws.Cells.Item(13, 13).Value = 147.91433
$ws.Cells.Item(13, 14).Value = 443.74299
$ws.Cells.Item(13, 15).Value = 0.2768624053389947
$ws.Cells.Item(13, 16).Value = 0.2935413991166814
$ws.Cells.Item(13, 17).Value = 4777.661258551359
$ws.Cells.Item(13, 18).Value = 42998.95132696224
$ws.Cells.Item(13, 19).Value = 0.0685574987163841
$ws.Cells.Item(13, 20).Value = 0.07454222340114451

$ws.Cells.Item(14, 7).Value = 32.300192
$ws.Cells.Item(14, 8).Value = 96.900576
$ws.Cells.Item(14, 9).Value = 0.2476229975407503
$ws.Cells.Item(14, 10).Value = 0.2539410918713864
$ws.Cells.Item(14, 13).Value = 83.50496933333334
$ws.Cells.Item(14, 14).Value = 250.514908
$ws.Cells.Item(14, 15).Value = 0.1563025480180701
$ws.Cells.Item(14, 16).Value = 0.1657186665504434
$ws.Cells.Item(14, 17).Value = 2697.226542420779
$ws.Cells.Item(14, 18).Value = 24275.03888178701
$ws.Cells.Item(14, 19).Value = 0.03870410546349157
$ws.Cells.Item(14, 20).Value = 0.04208277912728979

$ws.Cells.Item(15, 7).Value = 32.300192
$ws.Cells.Item(15, 8).Value = 96.900576
$ws.Cells.Item(15, 9).Value = 0.2476229975407503
$ws.Cells.Item(15, 10).Value = 0.2539410918713864
$ws.Cells.Item(15, 13).Value = 91.06846250000001
$ws.Cells.Item(15, 14).Value = 182.136925
$ws.Cells.Item(15, 15).Value = 0.1704597085236707
$ws.Cells.Item(15, 16).Value = 0.1204857969594293
$ws.Cells.Item(15, 17).Value = 2941.528823894801
$ws.Cells.Item(15, 18).Value = 17649.1729433688
$ws.Cells.Item(15, 19).Value = 0.04220974398455393
$ws.Cells.Item(15, 20).Value = 0.03059629483487165

$ws.Cells.Item(16, 7).Value = 32.300192
$ws.Cells.Item(16, 8).Value = 96.900576
$ws.Cells.Item(16, 9).Value = 0.2476229975407503
$ws.Cells.Item(16, 10).Value = 0.2539410918713864
$ws.Cells.Item(16, 13).Value = 89.83562999999999
$ws.Cells.Item(16, 14).Value = 269.50689
$ws.Cells.Item(16, 15).Value = 0.1681521229683693
$ws.Cells.Item(16, 16).Value = 0.1782820942415013
$ws.Cells.Item(16, 17).Value = 2901.70809744096
$ws.Cells.Item(16, 18).Value = 26115.37287696864
$ws.Cells.Item(16, 19).Value = 0.04163833273226845
$ws.Cells.Item(16, 20).Value = 0.04527314967280425

$ws.Cells.Item(17, 7).Value = 9.736177999999999
$ws.Cells.Item(17, 8).Value = 19.472356
$ws.Cells.Item(17, 9).Value = 0.07464047213559308
$ws.Cells.Item(17, 10).Value = 0.0510299478916239
$ws.Cells.Item(17, 13).Value = 121.928739
$ws.Cells.Item(17, 14).Value = 365.786217
$ws.Cells.Item(17, 15).Value = 0.2282232151508951
$ws.Cells.Item(17, 16).Value = 0.2419720431319445
$ws.Cells.Item(17, 17).Value = 1187.119906219542
$ws.Cells.Item(17, 18).Value = 7122.719437317251
$ws.Cells.Item(17, 19).Value = 0.01703468853116585
$ws.Cells.Item(17, 20).Value = 0.0123478207522529

$ws.Cells.Item(18, 7).Value = 9.736177999999999
$ws.Cells.Item(18, 8).Value = 19.472356
$ws.Cells.Item(18, 9).Value = 0.07464047213559308
$ws.Cells.Item(18, 10).Value = 0.0510299478916239
$ws.Cells.Item(18, 13).Value = 147.91433
$ws.Cells.Item(18, 14).Value = 443.74299
$ws.Cells.Item(18, 15).Value = 0.2768624053389947
$ws.Cells.Item(18, 16).Value = 0.2935413991166814
$ws.Cells.Item(18, 17).Value = 1440.12024563074
$ws.Cells.Item(18, 18).Value = 8640.721473784439
$ws.Cells.Item(18, 19).Value = 0.02066514065109851
$ws.Cells.Item(18, 20).Value = 0.01497940230095863

$ws.Cells.Item(19, 7).Value = 9.736177999999999
$ws.Cells.Item(19, 8).Value = 19.472356
$ws.Cells.Item(19, 9).Value = 0.07464047213559308
$ws.Cells.Item(19, 10).Value = 0.0510299478916239
$ws.Cells.Item(19, 13).Value = 83.50496933333334
$ws.Cells.Item(19, 14).Value = 250.514908
$ws.Cells.Item(19, 15).Value = 0.1563025480180701
$ws.Cells.Item(19, 16).Value = 0.1657186665504434
$ws.Cells.Item(19, 17).Value = 813.0192453138745
$ws.Cells.Item(19, 18).Value = 4878.115471883248
$ws.Cells.Item(19, 19).Value = 0.01166649598006496
$ws.Cells.Item(19, 20).Value = 0.008456614918738522

$ws.Cells.Item(20, 7).Value = 9.736177999999999
$ws.Cells.Item(20, 8).Value = 19.472356
$ws.Cells.Item(20, 9).Value = 0.07464047213559308
$ws.Cells.Item(20, 10).Value = 0.0510299478916239
$ws.Cells.Item(20, 13).Value = 91.06846250000001
$ws.Cells.Item(20, 14).Value = 182.136925
$ws.Cells.Item(20, 15).Value = 0.1704597085236707
$ws.Cells.Item(20, 16).Value = 0.1204857969594293
$ws.Cells.Item(20, 17).Value = 886.6587610863249
$ws.Cells.Item(20, 18).Value = 3546.6350443453
$ws.Cells.Item(20, 19).Value = 0.01272319312430236
$ws.Cells.Item(20, 20).Value = 0.006148383940520456

$ws.Cells.Item(21, 7).Value = 9.736177999999999
$ws.Cells.Item(21, 8).Value = 19.472356
$ws.Cells.Item(21, 9).Value = 0.07464047213559308
$ws.Cells.Item(21, 10).Value = 0.0510299478916239
$ws.Cells.Item(21, 13).Value = 89.83562999999999
$ws.Cells.Item(21, 14).Value = 269.50689
$ws.Cells.Item(21, 15).Value = 0.1681521229683693
$ws.Cells.Item(21, 16).Value = 0.1782820942415013
$ws.Cells.Item(21, 17).Value = 874.6556844221399
$ws.Cells.Item(21, 18).Value = 5247.934106532839
$ws.Cells.Item(21, 19).Value = 0.01255095384896139
$ws.Cells.Item(21, 20).Value = 0.009097725979153393

$ws.Cells.Item(22, 7).Value = 36.79199966666667
$ws.Cells.Item(22, 8).Value = 110.375999
$ws.Cells.Item(22, 9).Value = 0.2820585476079611
$ws.Cells.Item(22, 10).Value = 0.2892552641013719
$ws.Cells.Item(22, 13).Value = 121.928739
$ws.Cells.Item(22, 14).Value = 365.786217
$ws.Cells.Item(22, 15).Value = 0.2282232151508951
$ws.Cells.Item(22, 16).Value = 0.2419720431319445
$ws.Cells.Item(22, 17).Value = 4486.002124645087
$ws.Cells.Item(22, 18).Value = 40374.01912180578
$ws.Cells.Item(22, 19).Value = 0.0643723085958807
$ws.Cells.Item(22, 20).Value = 0.06999168724127916

$ws.Cells.Item(23, 7).Value = 36.79199966666667
$ws.Cells.Item(23, 8).Value = 110.375999
$ws.Cells.Item(23, 9).Value = 0.2820585476079611
$ws.Cells.Item(23, 10).Value = 0.2892552641013719
$ws.Cells.Item(23, 13).Value = 147.91433
$ws.Cells.Item(23, 14).Value = 443.74299
$ws.Cells.Item(23, 15).Value = 0.2768624053389947
$ws.Cells.Item(23, 16).Value = 0.2935413991166814
$ws.Cells.Item(23, 17).Value = 5442.063980055223
$ws.Cells.Item(23, 18).Value = 48978.57582049701
$ws.Cells.Item(23, 19).Value = 0.07809140793716346
$ws.Cells.Item(23, 20).Value = 0.08490839492618189

$ws.Cells.Item(24, 7).Value = 36.79199966666667
$ws.Cells.Item(24, 8).Value = 110.375999
$ws.Cells.Item(24, 9).Value = 0.2820585476079611
$ws.Cells.Item(24, 10).Value = 0.2892552641013719
$ws.Cells.Item(24, 13).Value = 83.50496933333334
$ws.Cells.Item(24, 14).Value = 250.514908
$ws.Cells.Item(24, 15).Value = 0.1563025480180701
$ws.Cells.Item(24, 16).Value = 0.1657186665504434
$ws.Cells.Item(24, 17).Value = 3072.31480387701
$ws.Cells.Item(24, 18).Value = 27650.83323489309
$ws.Cells.Item(24, 19).Value = 0.04408646968140044
$ws.Cells.Item(24, 20).Value = 0.04793499665957568

$ws.Cells.Item(25, 7).Value = 36.79199966666667
$ws.Cells.Item(25, 8).Value = 110.375999
$ws.Cells.Item(25, 9).Value = 0.2820585476079611
$ws.Cells.Item(25, 10).Value = 0.2892552641013719
$ws.Cells.Item(25, 13).Value = 91.06846250000001
$ws.Cells.Item(25, 14).Value = 182.136925
$ws.Cells.Item(25, 15).Value = 0.1704597085236707
$ws.Cells.Item(25, 16).Value = 0.1204857969594293
$ws.Cells.Item(25, 17).Value = 3350.590841943846
$ws.Cells.Item(25, 18).Value = 20103.54505166308
$ws.Cells.Item(25, 19).Value = 0.04807961781186296
$ws.Cells.Item(25, 20).Value = 0.034851151019964

$ws.Cells.Item(26, 7).Value = 36.79199966666667
$ws.Cells.Item(26, 8).Value = 110.375999
$ws.Cells.Item(26, 9).Value = 0.2820585476079611
$ws.Cells.Item(26, 10).Value = 0.2892552641013719
$ws.Cells.Item(26, 13).Value = 89.83562999999999
$ws.Cells.Item(26, 14).Value = 269.50689
$ws.Cells.Item(26, 15).Value = 0.1681521229683693
$ws.Cells.Item(26, 16).Value = 0.1782820942415013
$ws.Cells.Item(26, 17).Value = 3305.23246901479
$ws.Cells.Item(26, 18).Value = 29747.09222113311
$ws.Cells.Item(26, 19).Value = 0.04742874358165352
$ws.Cells.Item(26, 20).Value = 0.05156903425437113
